$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 with the same ticker value as row 2 ("GRT-USD"),
# extending the data range from A1:A2 to A1:A3.
$ws.Range("A3").Value = "GRT-USD"
